$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 33 with new Mac-Address / Document type data
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 10032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"
$ws.Range("G33").Value = "now()"

# Update selection to match post-edit state
$ws.Range("B30").Select()
